# Error Calculations and Plots
#
# The "missing data" mask for this combination/seed was regenerated:
#  - a handful of individual cells in column D (and F2) flip between
#    present/blank,
#  - row 26 ("RM 232") is dropped entirely, so every row below it
#    (previously 27-35) shifts up by one and the used range shrinks
#    from A1:F35 to A1:F33.
#
# We apply this by first collapsing the two now-unused trailing rows,
# then writing the final values for every row whose content differs
# from the original file (rows 26-33 are fully rewritten because the
# row-26 deletion shifts their identities; a few scattered cells above
# that just toggle blank/filled).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is now 2 rows shorter (A1:F35 -> A1:F33); drop the trailing
# rows that no longer exist.
$ws.Range("A34:F35").Delete()

# Scattered single-cell "missing value" toggles in rows 2-25.
$ws.Range("F2").ClearContents()
$ws.Range("D6").Value = -14.2
$ws.Range("D8").ClearContents()
$ws.Range("D18").Value = -15.2
$ws.Range("D20").ClearContents()
$ws.Range("D23").Value = -13.9
$ws.Range("D25").ClearContents()

# Row 26 ("RM 232") is gone, so rows 26-33 now hold what used to be
# rows 27-35 (with their own missing-value mask also changed).
$ws.Range("A26").Value = "SC 5"
$ws.Range("B26").Value = -20.2
$ws.Range("C26").Value = 10.8
$ws.Range("D26").Value = -13.8
$ws.Range("E26").Value = -5
$ws.Range("F26").Value = 17.38

$ws.Range("A27").Value = "SC 101"
$ws.Range("B27").Value = -20.4
$ws.Range("C27").Value = 10
$ws.Range("D27").Value = -14.6
$ws.Range("E27").Value = -10
$ws.Range("F27").Value = 17

$ws.Range("A28").Value = "SC 105"
$ws.Range("B28").ClearContents()
$ws.Range("C28").Value = 11.1
$ws.Range("D28").Value = -13.7
$ws.Range("E28").Value = -5.9
$ws.Range("F28").Value = 17.44

$ws.Range("A29").Value = "SC 119"
$ws.Range("B29").ClearContents()
$ws.Range("C29").Value = 11.2
$ws.Range("D29").Value = -13
$ws.Range("E29").Value = -6.8
$ws.Range("F29").Value = 18.06

$ws.Range("A30").Value = "SC 120"
$ws.Range("B30").Value = -19.7
$ws.Range("C30").Value = 11.4
$ws.Range("D30").Value = -13.6
$ws.Range("E30").Value = -5.7
$ws.Range("F30").Value = 16.89

$ws.Range("A31").Value = "SC 132"
$ws.Range("B31").Value = -18.8
$ws.Range("C31").Value = 15.3
$ws.Range("D31").Value = -13.7
$ws.Range("E31").Value = -8.1
$ws.Range("F31").Value = 17.18

$ws.Range("A32").Value = "SC 193"
$ws.Range("B32").ClearContents()
$ws.Range("C32").Value = 10.5
$ws.Range("D32").Value = -14.7
$ws.Range("E32").Value = -6.4
$ws.Range("F32").Value = 17.39

$ws.Range("A33").Value = "SC 232"
$ws.Range("B33").Value = -19.5
$ws.Range("C33").Value = 10.4
$ws.Range("D33").Value = -14.1
$ws.Range("E33").Value = -10.7
$ws.Range("F33").Value = 17.53
